# Update the Cort-Sstr2 LR-pairs sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new data only has 8 data rows (rows 2-9); the old trailing rows 10-11 are removed.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cort"
$ws.Range("C2").Value = "Sstr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2586766666666667
$ws.Range("H2").Value = 0.77603
$ws.Range("I2").Value = 0.9158046909573684
$ws.Range("J2").Value = 0.9158046909573684
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.193227
$ws.Range("N2").Value = 0.579681
$ws.Range("O2").Value = 0.1203927045071173
$ws.Range("P2").Value = 0.1203927045071173
$ws.Range("Q2").Value = 0.04998331627
$ws.Range("R2").Value = 0.44984984643
$ws.Range("S2").Value = 0.1102562035446623
$ws.Range("T2").Value = 0.1102562035446623

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cort"
$ws.Range("C3").Value = "Sstr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2586766666666667
$ws.Range("H3").Value = 0.77603
$ws.Range("I3").Value = 0.9158046909573684
$ws.Range("J3").Value = 0.9158046909573684
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.093643
$ws.Range("N3").Value = 3.280929
$ws.Range("O3").Value = 0.6814091122631787
$ws.Range("P3").Value = 0.6814091122631788
$ws.Range("Q3").Value = 0.2828999257633333
$ws.Range("R3").Value = 2.54609933187
$ws.Range("S3").Value = 0.6240376614717151
$ws.Range("T3").Value = 0.6240376614717152

# Row 4: FAPs -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cort"
$ws.Range("C4").Value = "Sstr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2586766666666667
$ws.Range("H4").Value = 0.77603
$ws.Range("I4").Value = 0.9158046909573684
$ws.Range("J4").Value = 0.9158046909573684
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2121483333333333
$ws.Range("N4").Value = 0.6364449999999999
$ws.Range("O4").Value = 0.132181898009478
$ws.Range("P4").Value = 0.132181898009478
$ws.Range("Q4").Value = 0.05487782370555554
$ws.Range("R4").Value = 0.49390041335
$ws.Range("S4").Value = 0.1210528022567284
$ws.Range("T4").Value = 0.1210528022567284

# Row 5: FAPs -> Resolving-Mac
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cort"
$ws.Range("C5").Value = "Sstr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2586766666666667
$ws.Range("H5").Value = 0.77603
$ws.Range("I5").Value = 0.9158046909573684
$ws.Range("J5").Value = 0.9158046909573684
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1059543333333333
$ws.Range("N5").Value = 0.317863
$ws.Range("O5").Value = 0.06601628522022597
$ws.Range("P5").Value = 0.06601628522022598
$ws.Range("Q5").Value = 0.02740791376555556
$ws.Range("R5").Value = 0.24667122389
$ws.Range("S5").Value = 0.06045802368426253
$ws.Range("T5").Value = 0.06045802368426254

# Row 6: Resolving-Mac -> ECs
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Cort"
$ws.Range("C6").Value = "Sstr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02378166666666667
$ws.Range("H6").Value = 0.071345
$ws.Range("I6").Value = 0.08419530904263167
$ws.Range("J6").Value = 0.08419530904263166
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.193227
$ws.Range("N6").Value = 0.579681
$ws.Range("O6").Value = 0.1203927045071173
$ws.Range("P6").Value = 0.1203927045071173
$ws.Range("Q6").Value = 0.004595260105000001
$ws.Range("R6").Value = 0.041357340945
$ws.Range("S6").Value = 0.01013650096245497
$ws.Range("T6").Value = 0.01013650096245497

# Row 7: Resolving-Mac -> FAPs
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Cort"
$ws.Range("C7").Value = "Sstr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02378166666666667
$ws.Range("H7").Value = 0.071345
$ws.Range("I7").Value = 0.08419530904263167
$ws.Range("J7").Value = 0.08419530904263166
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.093643
$ws.Range("N7").Value = 3.280929
$ws.Range("O7").Value = 0.6814091122631787
$ws.Range("P7").Value = 0.6814091122631788
$ws.Range("Q7").Value = 0.02600865327833334
$ws.Range("R7").Value = 0.234077879505
$ws.Range("S7").Value = 0.05737145079146363
$ws.Range("T7").Value = 0.05737145079146363

# Row 8: Resolving-Mac -> MuSCs
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cort"
$ws.Range("C8").Value = "Sstr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.02378166666666667
$ws.Range("H8").Value = 0.071345
$ws.Range("I8").Value = 0.08419530904263167
$ws.Range("J8").Value = 0.08419530904263166
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2121483333333333
$ws.Range("N8").Value = 0.6364449999999999
$ws.Range("O8").Value = 0.132181898009478
$ws.Range("P8").Value = 0.132181898009478
$ws.Range("Q8").Value = 0.005045240947222222
$ws.Range("R8").Value = 0.045407168525
$ws.Range("S8").Value = 0.01112909575274962
$ws.Range("T8").Value = 0.01112909575274962

# Row 9: Resolving-Mac -> Resolving-Mac
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cort"
$ws.Range("C9").Value = "Sstr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.02378166666666667
$ws.Range("H9").Value = 0.071345
$ws.Range("I9").Value = 0.08419530904263167
$ws.Range("J9").Value = 0.08419530904263166
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1059543333333333
$ws.Range("N9").Value = 0.317863
$ws.Range("O9").Value = 0.06601628522022597
$ws.Range("P9").Value = 0.06601628522022598
$ws.Range("Q9").Value = 0.002519770637222223
$ws.Range("R9").Value = 0.022677935735
$ws.Range("S9").Value = 0.005558261535963443
$ws.Range("T9").Value = 0.005558261535963444
